# Add "2022-Q3" quarterly data:
#  - insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q2")
#  - populate it with the fund-holdings table for that quarter
#  - insert a matching summary row at the top of the "总计" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet before the existing "2022-Q2" sheet
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)          # currently "2022-Q2"
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Borrow the cell formatting (fonts/borders/alignment) from the "2022-Q1"
# sheet, which already has the same 12-row (1 header + 11 data) shape, so
# the new sheet ends up with identical styling without inventing new xfs.
# (Copied in pieces so we don't manufacture an unused styled cell at A1,
# which the source sheet never populates either.)
$fmtSource = $wb.Worksheets.Item(4)            # "2022-Q1" (A1:H12)
$fmtSource.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)         # xlPasteFormats
$fmtSource.Range("A2:A12").Copy()
$q3.Range("A2:A12").PasteSpecial(-4122)
$fmtSource.Range("B2:G12").Copy()
$q3.Range("B2:G12").PasteSpecial(-4122)

# Columns B..G hold text-like values (fund code keeps leading zeros, the
# numeric-looking figures are stored as text, matching every other
# quarter sheet in this workbook) while A and H stay numeric.
$q3.Range("B2:G12").NumberFormat = "@"

# -- header row ---------------------------------------------------------
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# -- data rows ------------------------------------------------------------
$q3data = @(
  @(0,  "012930", "中庚价值先锋股票",             "47.83", "94.71", "4.24", "2.0280", 8),
  @(1,  "005416", "鹏华尊惠18个月定期开放混合A",  "7.53",  "37.52", "1.50", "0.1130", 6),
  @(2,  "003165", "鹏华弘嘉灵活配置混合A",        "0.82",  "91.09", "4.05", "0.0332", 2),
  @(3,  "009667", "鹏华安庆混合A",                "2.12",  "38.12", "1.17", "0.0248", 9),
  @(4,  "011572", "鹏华安荣混合A",                "1.50",  "39.92", "1.62", "0.0243", 5),
  @(5,  "009230", "鹏华安和混合A",                "1.80",  "39.30", "1.22", "0.0220", 9),
  @(6,  "009668", "鹏华安庆混合C",                "1.32",  "38.12", "1.17", "0.0154", 9),
  @(7,  "003166", "鹏华弘嘉灵活配置混合C",        "0.25",  "91.09", "4.05", "0.0101", 2),
  @(8,  "009231", "鹏华安和混合C",                "0.73",  "39.30", "1.22", "0.0089", 9),
  @(9,  "005417", "鹏华尊惠18个月定期开放混合C",  "0.58",  "37.52", "1.50", "0.0087", 6),
  @(10, "011573", "鹏华安荣混合C",                "0.34",  "39.92", "1.62", "0.0055", 5)
)

for ($i = 0; $i -lt $q3data.Length; $i++) {
  $r = $i + 2
  $rec = $q3data[$i]
  $q3.Cells.Item($r,1).Value = $rec[0]
  $q3.Cells.Item($r,2).Value = $rec[1]
  $q3.Cells.Item($r,3).Value = $rec[2]
  $q3.Cells.Item($r,4).Value = $rec[3]
  $q3.Cells.Item($r,5).Value = $rec[4]
  $q3.Cells.Item($r,6).Value = $rec[5]
  $q3.Cells.Item($r,7).Value = $rec[6]
  $q3.Cells.Item($r,8).Value = $rec[7]
}

# ---------------------------------------------------------------------------
# 2) Prepend the "2022-Q3" summary row on the "总计" sheet, shifting the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Preserve the "A" column's number style (bold/border) on the row that
# becomes newly used (old row 8 -> row 9) by copying an existing styled
# cell's format down before the values are rewritten.
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)         # xlPasteFormats

$totalData = @(
  @(0, "2022-Q3", 11, 2.29),
  @(1, "2022-Q2", 9,  2.71),
  @(2, "2022-Q1", 11, 3.66),
  @(3, "2021-Q4", 20, 5.22),
  @(4, "2021-Q3", 6,  1.08),
  @(5, "2021-Q2", 11, 1.45),
  @(6, "2021-Q1", 6,  1.42),
  @(7, "2020-Q4", 7,  3.87)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
  $r = $i + 2
  $rec = $totalData[$i]
  $total.Cells.Item($r,1).Value = $rec[0]
  $total.Cells.Item($r,2).Value = $rec[1]
  $total.Cells.Item($r,3).Value = $rec[2]
  $total.Cells.Item($r,4).Value = $rec[3]
}

Write-Output "2022-Q3 sheet + summary row added"
